$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1101.9286
$ws.Range("I28").Value = 270.16666
$ws.Range("J28").Value = 1725.75
$ws.Range("K28").Value = 270.16666
$ws.Range("L28").Value = 1725.75
$ws.Range("M28").Value = 214.83334
$ws.Range("N28").Value = -2695.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 2800
$ws.Range("J29").Value = 2800
$ws.Range("L29").Value = 8400
$ws.Range("N29").Value = -8962

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 679.8
$ws.Range("I38").Value = 349.75
$ws.Range("J38").Value = 2000
$ws.Range("K38").Value = 1049.25
$ws.Range("L38").Value = 6000
$ws.Range("M38").Value = -677.25
$ws.Range("N38").Value = -6744

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 57693360
$ws.Range("I135").Value = 23810668
$ws.Range("J135").Value = 200000660
$ws.Range("K135").Value = 214296012
$ws.Range("L135").Value = 1800005940
$ws.Range("M135").Value = -214293477
$ws.Range("N135").Value = -1800011010

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2189.2327
$ws.Range("I137").Value = 1824.1538
$ws.Range("J137").Value = 2747.5881
$ws.Range("K137").Value = 5472.4614
$ws.Range("L137").Value = 8242.764299999999
$ws.Range("M137").Value = -2922.4614
$ws.Range("N137").Value = -13342.7643

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2558.5
$ws.Range("I63").Value = 2731
$ws.Range("J63").Value = 1006
$ws.Range("K63").Value = 2731
$ws.Range("L63").Value = 1006
$ws.Range("M63").Value = -2045
$ws.Range("N63").Value = -2378

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2558.5
$ws.Range("I66").Value = 2731
$ws.Range("J66").Value = 1006
$ws.Range("K66").Value = 13655
$ws.Range("L66").Value = 5030
$ws.Range("M66").Value = -10223
$ws.Range("N66").Value = -11894

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 7631.75
$ws.Range("I74").Value = 2769.5264
$ws.Range("J74").Value = 100014
$ws.Range("K74").Value = 2769.5264
$ws.Range("L74").Value = 100014
$ws.Range("M74").Value = -1895.5264
$ws.Range("N74").Value = -101762

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 7631.75
$ws.Range("I77").Value = 2769.5264
$ws.Range("J77").Value = 100014
$ws.Range("K77").Value = 13847.632
$ws.Range("L77").Value = 500070
$ws.Range("M77").Value = -9479.632000000001
$ws.Range("N77").Value = -508806

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 36073.332
$ws.Range("J35").Value = 36073.332
$ws.Range("L35").Value = 36073.332
$ws.Range("N35").Value = -36693.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 21580.04
$ws.Range("I134").Value = 1713.15
$ws.Range("J134").Value = 93823.27
$ws.Range("K134").Value = 5139.450000000001
$ws.Range("L134").Value = 281469.81
$ws.Range("M134").Value = -2604.450000000001
$ws.Range("N134").Value = -286539.81

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2311.7908
$ws.Range("I31").Value = 1689.0541
$ws.Range("J31").Value = 6152
$ws.Range("K31").Value = 1689.0541
$ws.Range("L31").Value = 6152
$ws.Range("M31").Value = -1394.0541
$ws.Range("N31").Value = -6742

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2311.7908
$ws.Range("I34").Value = 1689.0541
$ws.Range("J34").Value = 6152
$ws.Range("K34").Value = 1689.0541
$ws.Range("L34").Value = 6152
$ws.Range("M34").Value = -1487.0541
$ws.Range("N34").Value = -6556

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2839.5557
$ws.Range("I62").Value = 2931.25
$ws.Range("J62").Value = 2106
$ws.Range("K62").Value = 2931.25
$ws.Range("L62").Value = 2106
$ws.Range("M62").Value = -2307.25
$ws.Range("N62").Value = -3354

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 2839.5557
$ws.Range("I65").Value = 2931.25
$ws.Range("J65").Value = 2106
$ws.Range("K65").Value = 14656.25
$ws.Range("L65").Value = 10530
$ws.Range("M65").Value = -11536.25
$ws.Range("N65").Value = -16770

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3896.7544
$ws.Range("I132").Value = 4318.5674
$ws.Range("J132").Value = 3116.4
$ws.Range("K132").Value = 12955.7022
$ws.Range("L132").Value = 9349.200000000001
$ws.Range("M132").Value = -10425.7022
$ws.Range("N132").Value = -14409.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2449.242
$ws.Range("I134").Value = 1459.027
$ws.Range("J134").Value = 3914.76
$ws.Range("K134").Value = 4377.081
$ws.Range("L134").Value = 11744.28
$ws.Range("M134").Value = -1842.081
$ws.Range("N134").Value = -16814.28

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 2140.9375
$ws.Range("I118").Value = 806.1111
$ws.Range("J118").Value = 3857.1428
$ws.Range("K118").Value = 2418.3333
$ws.Range("L118").Value = 11571.4284
$ws.Range("M118").Value = -1175.3333
$ws.Range("N118").Value = -14057.4284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 30
$ws.Range("I121").Value = 30
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 90
$ws.Range("L121").Value = 0
$ws.Range("M121").Value = 1220
$ws.Range("N121").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 3004.2222
$ws.Range("I123").Value = 1956
$ws.Range("J123").Value = 3242.4546
$ws.Range("K123").Value = 5868
$ws.Range("L123").Value = 9727.363799999999
$ws.Range("M123").Value = -3418
$ws.Range("N123").Value = -14627.3638

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 2322
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 2322
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 6966
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -16806

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3307
$ws.Range("I102").Value = 2912.625
$ws.Range("J102").Value = 4167.4546
$ws.Range("K102").Value = 2912.625
$ws.Range("L102").Value = 4167.4546
$ws.Range("M102").Value = -1290.625
$ws.Range("N102").Value = -7411.4546

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 11025.5
$ws.Range("I122").Value = 18667
$ws.Range("J122").Value = 3384
$ws.Range("K122").Value = 56001
$ws.Range("L122").Value = 10152
$ws.Range("M122").Value = -53551
$ws.Range("N122").Value = -15052

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 6506.3335
$ws.Range("I132").Value = 4419.8423
$ws.Range("J132").Value = 26328
$ws.Range("K132").Value = 13259.5269
$ws.Range("L132").Value = 78984
$ws.Range("M132").Value = -10729.5269
$ws.Range("N132").Value = -84044

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1618.625
$ws.Range("I132").Value = 744.11536
$ws.Range("J132").Value = 3242.7144
$ws.Range("K132").Value = 2232.34608
$ws.Range("L132").Value = 9728.143199999999
$ws.Range("M132").Value = 297.6539199999997
$ws.Range("N132").Value = -14788.1432
